$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$ws.Range("A3").Value = 108127462
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
